# Trade #15 closed at 2026-02-16 21:24:21 - leadlag DOWN +0.000%
# Appends a new trade row (row 14) to the "leadlag" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("leadlag")

$row = 14

# Trade #
$ws.Cells.Item($row, 1).Value = 15

# Date / Time - force text so the "yyyy-mm-dd" / "hh:mm:ss" strings are not
# auto-converted into date/time serial numbers. ClearFormats() afterwards
# drops the temporary "@" number-format style so the cell keeps the sheet's
# default (un-styled) look, matching the other rows.
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2026-02-16"
$ws.Cells.Item($row, 2).ClearFormats()
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "21:24:21"
$ws.Cells.Item($row, 3).ClearFormats()

# Strategy / Side
$ws.Cells.Item($row, 4).Value = "leadlag"
$ws.Cells.Item($row, 5).Value = "DOWN"

# Entry Price
$ws.Cells.Item($row, 6).Value = 69347.99000000001

# Exit Price (blank/open trade)
$ws.Cells.Item($row, 7).Value = ""

# Status
$ws.Cells.Item($row, 8).Value = "OPEN"

# P&L % / P&L $
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0

# Confidence
$ws.Cells.Item($row, 11).Value = 0.6074000000000001

# Entry Reason
$ws.Cells.Item($row, 12).Value = "Binance leading with -0.061% move"

# Exit Reason (blank/open trade)
$ws.Cells.Item($row, 13).Value = ""

# Duration (min)
$ws.Cells.Item($row, 14).Value = 0
